$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 180, shifting rows 180-222 down to 181-223.
$ws.Rows(180).Insert()

# Populate the newly inserted row 180 with the new record.
$ws.Range("A180").Value = 5
$ws.Range("B180").Value = "Macroferia Regional de Talca"
$ws.Range("C180").Value = "Maule"
$ws.Range("D180").Value = 44511
$ws.Range("E180").Value = 7
$ws.Range("F180").Value = 100112023
$ws.Range("G180").Value = "Brócoli"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 5000
$ws.Range("K180").Value = 600
$ws.Range("L180").Value = 600
$ws.Range("M180").Value = 600
$ws.Range("N180").Value = "$/unidad"
$ws.Range("O180").Value = "Región del Maule"
$ws.Range("P180").Value = 600
$ws.Range("Q180").Value = 1
$ws.Range("R180").Value = "Hortaliza"
